# Auto commit at 2025-08-31  8:09:11.14
# Append two new daily records (row 60: 四方坪站, row 61: 高岭站) for 2025-08-30
# to the "daydata" sheet, then move the active selection to H55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60 - 四方坪站 (Sifangping station)
$ws.Cells.Item(60, 1).Value = 45899
$ws.Cells.Item(60, 2).Value = "四方坪站"
$ws.Cells.Item(60, 3).Value = 11907.25
$ws.Cells.Item(60, 4).Value = 9954.78
$ws.Cells.Item(60, 5).Value = 4153.97
$ws.Cells.Item(60, 6).Value = 468

# Row 61 - 高岭站 (Gaoling station)
$ws.Cells.Item(61, 1).Value = 45899
$ws.Cells.Item(61, 2).Value = "高岭站"
$ws.Cells.Item(61, 3).Value = 6277.77
$ws.Cells.Item(61, 4).Value = 5238.43
$ws.Cells.Item(61, 5).Value = 1587.25
$ws.Cells.Item(61, 6).Value = 192

# Move the selection, as recorded in the saved view state.
$ws.Range("H55").Select() | Out-Null
